$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "About"
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Header: add state label "Oregon" next to the title, bump the "last updated" date
$about.Range("B1").Value = "Oregon"
$about.Range("C1").Value = 44834

# Insert 3 new rows before the old "time series" note (2 lines of new text +
# 1 blank separator row), pushing the old note down.
$about.Range("A55:A57").EntireRow.Insert()
$about.Range("A55").Value = "The units for rail in the BTS data set are unclear whether they report locomotive-miles or train-car-miles. Based on the"
$about.Range("A56").Value = "scale, we assume the units are train-car-miles."

# ---------------------------------------------------------------------------
# Sheet "BTS NTS Modal Profile Data"
# ---------------------------------------------------------------------------
$bts = $wb.Worksheets.Item("BTS NTS Modal Profile Data")

# Insert a new assumption row before the "Intercity (Amtrak)" section.
$bts.Rows.Item(22).Insert()
$bts.Range("A22").Value = "Assumption - train cars per locomotive"
$bts.Range("B22").Value = 10

# The old final "weighted value, adjusted for number of train cars per
# locomotive" row (now at row 38 after the insert above) is removed - its
# "/10" adjustment is folded into the weighted-average formula below instead,
# using the new assumption cell B22.
$bts.Rows.Item(38).Delete()

$bts.Range("B37").Formula = "=(B26*B25+B34*B28+B35*B29+B36*B30)/SUM(B25,B28:B30)*B22"

# ---------------------------------------------------------------------------
# Sheet "AVLo-passengers"
# ---------------------------------------------------------------------------
$pass = $wb.Worksheets.Item("AVLo-passengers")

# B37 on the BTS sheet no longer pre-divides by 10, so do it here instead.
$pass.Range("B5").Formula = "='BTS NTS Modal Profile Data'!B37/10"
